# Updating deaths and births excel file
# - Fill in the newly-available Births/Deaths values in column C for years
#   2012-2022 on the "SPAIN_Components_pop_change" sheet (matching the
#   existing formatting used by the sibling cells directly above, e.g. C22/C23).
# - Refresh the view state (zoom 95% -> 120%, selections) across all sheets.
# - Tweak the header/footer font style name and a stray row height on the
#   "Template comp pop change" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. New Births / Deaths figures on "SPAIN_Components_pop_change" (sheet 1)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$newValues = @{
    "C41" = 454648;  "C42" = 402950;
    "C50" = 425715;  "C51" = 390419;
    "C59" = 427595;  "C60" = 395830;
    "C68" = 420290;  "C69" = 422568;
    "C77" = 410583;  "C78" = 410611;
    "C86" = 393181;  "C87" = 424523;
    "C95" = 372777;  "C96" = 427721;
    "C104" = 360617; "C105" = 418703;
    "C113" = 341315; "C114" = 493776;
    "C122" = 337380; "C123" = 450744;
    "C131" = 329251; "C132" = 464417;
}

# C22/C23 already carry the formatting (font/fill/border/number format) that
# the newly populated cells should match, so copy it across before writing
# the values in place.
$formatSource = $ws1.Range("C22")

foreach ($addr in $newValues.Keys) {
    $formatSource.Copy()
    $ws1.Range($addr).PasteSpecial(-4122)
    $ws1.Range($addr).Value = $newValues[$addr]
}
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 2. View state refresh: zoom 95 -> 120 on every sheet, plus updated
#    scroll position / selection per sheet.
# ---------------------------------------------------------------------------

# Sheet 1: SPAIN_Components_pop_change
$ws1.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 115
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("C132").Select()

# Sheet 2: INE Total Population SPAIN
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws2.Range("B26").Select()

# Sheet 3: INE National Increase SPAIN
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("C19").Select()

# Sheet 4: INE Net External Migration
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws4.Range("I30").Select()

# Sheet 5: Template comp pop change
$ws5 = $wb.Worksheets.Item(5)
$ws5.Activate()
$excel.ActiveWindow.Zoom = 120
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws5.Range("C29").Select()

# Row 4 on this sheet no longer needs the tall, wrapped height.
$ws5.Rows.Item(4).RowHeight = 12.8

# ---------------------------------------------------------------------------
# 3. Header / footer font style name: "...,Normal" -> "...,Regular"
#    (applies to every sheet's odd header/footer).
# ---------------------------------------------------------------------------
foreach ($ws in @($ws1, $ws2, $ws3, $ws4, $ws5)) {
    $ws.PageSetup.CenterHeader = "&""Times New Roman,Regular""&12&A"
    $ws.PageSetup.CenterFooter = "&""Times New Roman,Regular""&12Página &P"
}

# Leave the originally active sheet/tab selected.
$ws1.Activate()
